# Apply the change described by the diff:
# - Split the "Programa" paragraph (PT) numbered items 1-5 with line breaks
# - Split the "Programa" paragraph (EN, italic) numbered items 1-5 with line breaks
# - Split the "Bibliografia" paragraph entries with double line breaks

$d = $word.ActiveDocument

function Replace-WithBreak {
    param(
        [string]$OldText,
        [string]$NewText
    )
    $r = $d.Content
    $find = $r.Find
    $result = $find.Execute($OldText, $true, $false, $false, $false, $false, $true, 1, $false, $NewText, 2)
    if (-not $result) {
        Write-Host "WARNING: replace failed for: $OldText"
    }
}

# --- Programa (Portuguese) ---
Replace-WithBreak "dinâmicas.2." "dinâmicas.^l2."
Replace-WithBreak "Pessoas3." "Pessoas^l3."
Replace-WithBreak "Marketing4." "Marketing^l4."
Replace-WithBreak "marketing)5." "marketing)^l5."

# --- Programa (English, italic) ---
Replace-WithBreak "capabilities.2." "capabilities.^l2."
Replace-WithBreak "Management3." "Management^l3."
Replace-WithBreak "basics4." "basics^l4."
Replace-WithBreak "plan)5." "plan)^l5."

# --- Bibliografia: insert a double line break between each reference entry ---
Replace-WithBreak ": Manole, 2014.Chiavenato, I. " ": Manole, 2014.^l^lChiavenato, I. "
Replace-WithBreak ", Campus, 2015.ROBBINS, S. P.;" ", Campus, 2015.^l^lROBBINS, S. P.;"
Replace-WithBreak " saraiva, 2013.KOTLER, P. - AR" " saraiva, 2013.^l^lKOTLER, P. - AR"
Replace-WithBreak " Pearson, 2014.KOTLER, P.; KEL" " Pearson, 2014.^l^lKOTLER, P.; KEL"
Replace-WithBreak " Pearson, 2019.CHIAVENATO, I. " " Pearson, 2019.^l^lCHIAVENATO, I. "
Replace-WithBreak " Manole, 2014. MAXIMIANO, A. C" " Manole, 2014. ^l^lMAXIMIANO, A. C"
Replace-WithBreak "o: Atlas, 2017.GUERRINI, F. M." "o: Atlas, 2017.^l^lGUERRINI, F. M."
Replace-WithBreak ": Campus, 2016.CHIAVENATO, I. " ": Campus, 2016.^l^lCHIAVENATO, I. "
Replace-WithBreak ": Manole, 2011.SILVA, M. M. L." ": Manole, 2011.^l^lSILVA, M. M. L."
Replace-WithBreak "Brasport, 2018.BOLMAN, L.G.; D" "Brasport, 2018.^l^lBOLMAN, L.G.; D"
Replace-WithBreak "ohn Wiley, 2013KOTLER, P.. O M" "ohn Wiley, 2013^l^lKOTLER, P.. O M"
Replace-WithBreak ". Bookman, 2005MINTZBERG, H. C" ". Bookman, 2005^l^lMINTZBERG, H. C"
